$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New crime data collected - weekly CompStat refresh:
#   Volume 30 Number 47 -> 48
#   Report week 11/20/2023-11/26/2023 -> 11/27/2023-12/3/2023
#   Updated weekly crime-complaint figures for rows 15-30
# ---------------------------------------------------------------------------

# --- Header text updates (volume number + reporting week dates) -----------
$volCell = $ws.Range("A8")
$volCell.Characters(21, 2).Text = "48"

$weekCell = $ws.Range("C9")
$weekCell.Characters(27, 10).Text = "11/27/2023"
$weekCell.Characters(48, 10).Text = "12/3/2023"

# --- Cells whose type/style changes (numeric <-> "N/A" text placeholder) --
$zeroSrc = $ws.Range("C23")
$naSrc = $ws.Range("E23")
$numSrc = $ws.Range("D27")

function Set-NAZero($ref) {
    $dst = $ws.Range($ref)
    $zeroSrc.Copy()
    $dst.PasteSpecial(-4122) | Out-Null
    $zeroSrc.Copy()
    $dst.PasteSpecial(-4163) | Out-Null
}

function Set-NADash($ref) {
    $dst = $ws.Range($ref)
    $naSrc.Copy()
    $dst.PasteSpecial(-4122) | Out-Null
    $naSrc.Copy()
    $dst.PasteSpecial(-4163) | Out-Null
}

function Set-Number($ref, $value) {
    $dst = $ws.Range($ref)
    $numSrc.Copy()
    $dst.PasteSpecial(-4122) | Out-Null
    $dst.Value = $value
}

Set-NAZero("D15")
Set-NADash("E15")

Set-Number("C22", 1)

Set-NAZero("C26")
Set-NAZero("D26")
Set-NADash("E26")

Set-Number("C27", 2)

$updates = @(
    @("G15", 1),
    @("C16", 4),
    @("E16", 0),
    @("G16", 19),
    @("H16", 5.263157894736),
    @("I16", 181),
    @("J16", 237),
    @("K16", -23.628691983122),
    @("L16", 6.470588235294),
    @("M16", 31.159420289855),
    @("N16", -77.403245942571),
    @("C17", 3),
    @("D17", 6),
    @("E17", -50),
    @("F17", 10),
    @("H17", -16.666666666666),
    @("I17", 184),
    @("J17", 173),
    @("K17", 6.358381502890),
    @("L17", 10.843373493975),
    @("M17", 97.849462365591),
    @("N17", -36.332179930795),
    @("C18", 2),
    @("D18", 11),
    @("E18", -81.818181818181),
    @("F18", 21),
    @("G18", 27),
    @("H18", -22.222222222222),
    @("I18", 270),
    @("J18", 406),
    @("K18", -33.497536945812),
    @("L18", 18.942731277533),
    @("M18", 56.976744186046),
    @("N18", -64.379947229551),
    @("C19", 22),
    @("D19", 29),
    @("E19", -24.137931034482),
    @("F19", 89),
    @("G19", 113),
    @("H19", -21.238938053097),
    @("I19", 1224),
    @("J19", 1285),
    @("K19", -4.747081712062),
    @("L19", 54.350567465321),
    @("M19", 25.153374233128),
    @("N19", -48.914858096828),
    @("G20", 2),
    @("H20", 0),
    @("I20", 45),
    @("K20", -18.181818181818),
    @("L20", -15.094339622641),
    @("M20", 9.756097560975),
    @("N20", -93.835616438356),
    @("C21", 32),
    @("D21", 50),
    @("E21", -36),
    @("F21", 142),
    @("G21", 174),
    @("H21", -18.390804597701),
    @("I21", 1913),
    @("J21", 2171),
    @("K21", -11.883924458774),
    @("L21", 34.813248766737),
    @("M21", 33.310104529616),
    @("N21", -61.632571199358),
    @("D22", 3),
    @("E22", -66.666666666666),
    @("F22", 2),
    @("G22", 5),
    @("H22", -60),
    @("I22", 40),
    @("J22", 41),
    @("K22", -2.439024390243),
    @("L22", 33.333333333333),
    @("M22", -11.111111111111),
    @("C24", 41),
    @("D24", 39),
    @("E24", 5.128205128205),
    @("F24", 146),
    @("G24", 154),
    @("H24", -5.194805194805),
    @("I24", 1941),
    @("J24", 1829),
    @("K24", 6.123564789502),
    @("L24", 61.212624584717),
    @("M24", 40.04329004329),
    @("C25", 10),
    @("D25", 8),
    @("E25", 25),
    @("F25", 32),
    @("G25", 34),
    @("H25", -5.882352941176),
    @("I25", 400),
    @("J25", 391),
    @("K25", 2.301790281329),
    @("L25", 37.457044673539),
    @("M25", 68.067226890756),
    @("G26", 1),
    @("H26", 0),
    @("L26", 14.285714285714),
    @("D27", 2),
    @("E27", 0),
    @("G27", 6),
    @("H27", -50),
    @("I27", 60),
    @("J27", 73),
    @("K27", -17.808219178082),
    @("L27", -7.692307692307),
    @("G30", 1),
    @("I30", 10),
    @("K30", -16.666666666666),
    @("L30", -28.571428571428)
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

Write-Output "done"
